$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap F:V content between row pairs that were reordered ---
$tmp = $ws.Range("F19:V19").Value2
$ws.Range("F19:V19").Value2 = $ws.Range("F20:V20").Value2
$ws.Range("F20:V20").Value2 = $tmp

$tmp = $ws.Range("F22:V22").Value2
$ws.Range("F22:V22").Value2 = $ws.Range("F23:V23").Value2
$ws.Range("F23:V23").Value2 = $tmp

$tmp = $ws.Range("F37:V37").Value2
$ws.Range("F37:V37").Value2 = $ws.Range("F38:V38").Value2
$ws.Range("F38:V38").Value2 = $tmp

$tmp = $ws.Range("F74:V74").Value2
$ws.Range("F74:V74").Value2 = $ws.Range("F75:V75").Value2
$ws.Range("F75:V75").Value2 = $tmp

$tmp = $ws.Range("F96:V96").Value2
$ws.Range("F96:V96").Value2 = $ws.Range("F97:V97").Value2
$ws.Range("F97:V97").Value2 = $tmp

$tmp = $ws.Range("F98:V98").Value2
$ws.Range("F98:V98").Value2 = $ws.Range("F99:V99").Value2
$ws.Range("F99:V99").Value2 = $tmp

$tmp = $ws.Range("F102:V102").Value2
$ws.Range("F102:V102").Value2 = $ws.Range("F103:V103").Value2
$ws.Range("F103:V103").Value2 = $tmp

$tmp = $ws.Range("F106:V106").Value2
$ws.Range("F106:V106").Value2 = $ws.Range("F107:V107").Value2
$ws.Range("F107:V107").Value2 = $tmp

$tmp = $ws.Range("F110:V110").Value2
$ws.Range("F110:V110").Value2 = $ws.Range("F111:V111").Value2
$ws.Range("F111:V111").Value2 = $tmp

$tmp = $ws.Range("F133:V133").Value2
$ws.Range("F133:V133").Value2 = $ws.Range("F134:V134").Value2
$ws.Range("F134:V134").Value2 = $tmp

$tmp = $ws.Range("F140:V140").Value2
$ws.Range("F140:V140").Value2 = $ws.Range("F141:V141").Value2
$ws.Range("F141:V141").Value2 = $tmp

$tmp = $ws.Range("F167:V167").Value2
$ws.Range("F167:V167").Value2 = $ws.Range("F168:V168").Value2
$ws.Range("F168:V168").Value2 = $tmp

$tmp = $ws.Range("F184:V184").Value2
$ws.Range("F184:V184").Value2 = $ws.Range("F186:V186").Value2
$ws.Range("F186:V186").Value2 = $tmp

$tmp = $ws.Range("F185:V185").Value2
$ws.Range("F185:V185").Value2 = $ws.Range("F187:V187").Value2
$ws.Range("F187:V187").Value2 = $tmp

$tmp = $ws.Range("F188:V188").Value2
$ws.Range("F188:V188").Value2 = $ws.Range("F191:V191").Value2
$ws.Range("F191:V191").Value2 = $tmp

$tmp = $ws.Range("F195:V195").Value2
$ws.Range("F195:V195").Value2 = $ws.Range("F196:V196").Value2
$ws.Range("F196:V196").Value2 = $tmp

# --- Step 2: append new rows 197-203, copying formatting (and the B/C/D text columns,
# which are identical on every row) from row 196 ---
$ws.Range("A196:V196").Copy($ws.Range("A197:V197"))
$ws.Range("A196:V196").Copy($ws.Range("A198:V198"))
$ws.Range("A196:V196").Copy($ws.Range("A199:V199"))
$ws.Range("A196:V196").Copy($ws.Range("A200:V200"))
$ws.Range("A196:V196").Copy($ws.Range("A201:V201"))
$ws.Range("A196:V196").Copy($ws.Range("A202:V202"))
$ws.Range("A196:V196").Copy($ws.Range("A203:V203"))

# row 197
$ws.Range("A197").Value2 = 196
$ws.Range("E197").Value2 = 45262.9375
$ws.Range("F197").Value2 = 'Huracan'
$ws.Range("G197").Value2 = 1
$ws.Range("H197").Value2 = 'Platense'
$ws.Range("I197").Value2 = 1
$ws.Range("J197").Value2 = 1.98
$ws.Range("K197").Value2 = '29/11/2023 03:43'
$ws.Range("L197").Value2 = 2.55
$ws.Range("M197").Value2 = '02/12/2023 22:29'
$ws.Range("N197").Value2 = 2.99
$ws.Range("O197").Value2 = '29/11/2023 03:43'
$ws.Range("P197").Value2 = 2.68
$ws.Range("Q197").Value2 = '02/12/2023 22:29'
$ws.Range("R197").Value2 = 4.88
$ws.Range("S197").Value2 = '29/11/2023 03:43'
$ws.Range("T197").Value2 = 3.73
$ws.Range("U197").Value2 = '02/12/2023 22:29'
$ws.Range("V197").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/huracan-platense/SzIUwiaU/'

# row 198
$ws.Range("A198").Value2 = 197
$ws.Range("E198").Value2 = 45263.0625
$ws.Range("F198").Value2 = 'Godoy Cruz'
$ws.Range("G198").Value2 = 0
$ws.Range("H198").Value2 = 'Banfield'
$ws.Range("I198").Value2 = 0
$ws.Range("J198").Value2 = 2.18
$ws.Range("K198").Value2 = '29/11/2023 03:43'
$ws.Range("L198").Value2 = 2.74
$ws.Range("M198").Value2 = '03/12/2023 01:29'
$ws.Range("N198").Value2 = 2.94
$ws.Range("O198").Value2 = '29/11/2023 03:43'
$ws.Range("P198").Value2 = 2.7
$ws.Range("Q198").Value2 = '03/12/2023 01:29'
$ws.Range("R198").Value2 = 4.08
$ws.Range("S198").Value2 = '29/11/2023 03:43'
$ws.Range("T198").Value2 = 3.36
$ws.Range("U198").Value2 = '03/12/2023 01:27'
$ws.Range("V198").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-banfield/lpJQvXqO/'

# row 199
$ws.Range("A199").Value2 = 198
$ws.Range("E199").Value2 = 45263.9375
$ws.Range("F199").Value2 = 'River Plate'
$ws.Range("G199").Value2 = 2
$ws.Range("H199").Value2 = 'Belgrano'
$ws.Range("I199").Value2 = 1
$ws.Range("J199").Value2 = 1.5
$ws.Range("K199").Value2 = '29/11/2023 03:43'
$ws.Range("L199").Value2 = 1.58
$ws.Range("M199").Value2 = '03/12/2023 22:12'
$ws.Range("N199").Value2 = 4.26
$ws.Range("O199").Value2 = '29/11/2023 03:43'
$ws.Range("P199").Value2 = 3.87
$ws.Range("Q199").Value2 = '03/12/2023 22:29'
$ws.Range("R199").Value2 = 6.07
$ws.Range("S199").Value2 = '29/11/2023 03:43'
$ws.Range("T199").Value2 = 6.73
$ws.Range("U199").Value2 = '03/12/2023 22:29'
$ws.Range("V199").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/river-plate-ca-belgrano-de-cordoba/fNJMuDUH/'

# row 200
$ws.Range("A200").Value2 = 199
$ws.Range("E200").Value2 = 45264.0625
$ws.Range("F200").Value2 = 'Racing Club'
$ws.Range("G200").Value2 = 2
$ws.Range("H200").Value2 = 'Rosario Central'
$ws.Range("I200").Value2 = 2
$ws.Range("J200").Value2 = 1.88
$ws.Range("K200").Value2 = '29/11/2023 03:43'
$ws.Range("L200").Value2 = 1.93
$ws.Range("M200").Value2 = '04/12/2023 01:29'
$ws.Range("N200").Value2 = 3.5
$ws.Range("O200").Value2 = '29/11/2023 03:43'
$ws.Range("P200").Value2 = 3.3
$ws.Range("Q200").Value2 = '04/12/2023 01:29'
$ws.Range("R200").Value2 = 4.03
$ws.Range("S200").Value2 = '29/11/2023 03:43'
$ws.Range("T200").Value2 = 4.62
$ws.Range("U200").Value2 = '04/12/2023 01:29'
$ws.Range("V200").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/racing-club-rosario-central/tKFItgFB/'

# row 201
$ws.Range("A201").Value2 = 200
$ws.Range("E201").Value2 = 45269.91666666666
$ws.Range("F201").Value2 = 'Godoy Cruz'
$ws.Range("G201").Value2 = 1
$ws.Range("H201").Value2 = 'Platense'
$ws.Range("I201").Value2 = 1
$ws.Range("J201").Value2 = 2.47
$ws.Range("K201").Value2 = '06/12/2023 13:11'
$ws.Range("L201").Value2 = 2.44
$ws.Range("M201").Value2 = '09/12/2023 21:27'
$ws.Range("N201").Value2 = 2.95
$ws.Range("O201").Value2 = '06/12/2023 13:11'
$ws.Range("P201").Value2 = 2.8
$ws.Range("Q201").Value2 = '09/12/2023 21:17'
$ws.Range("R201").Value2 = 3.16
$ws.Range("S201").Value2 = '06/12/2023 13:11'
$ws.Range("T201").Value2 = 3.75
$ws.Range("U201").Value2 = '09/12/2023 21:27'
$ws.Range("V201").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/godoy-cruz-platense/jcFQNmw2/'

# row 202
$ws.Range("A202").Value2 = 201
$ws.Range("E202").Value2 = 45270.08333333334
$ws.Range("F202").Value2 = 'River Plate'
$ws.Range("G202").Value2 = 0
$ws.Range("H202").Value2 = 'Rosario Central'
$ws.Range("I202").Value2 = 0
$ws.Range("J202").Value2 = 1.61
$ws.Range("K202").Value2 = '06/12/2023 13:12'
$ws.Range("L202").Value2 = 1.62
$ws.Range("M202").Value2 = '10/12/2023 01:59'
$ws.Range("N202").Value2 = 3.66
$ws.Range("O202").Value2 = '06/12/2023 13:12'
$ws.Range("P202").Value2 = 4.15
$ws.Range("Q202").Value2 = '10/12/2023 01:59'
$ws.Range("R202").Value2 = 5.81
$ws.Range("S202").Value2 = '06/12/2023 13:12'
$ws.Range("T202").Value2 = 5.59
$ws.Range("U202").Value2 = '10/12/2023 01:56'
$ws.Range("V202").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/river-plate-rosario-central/x0z1nXF2/'

# row 203
$ws.Range("A203").Value2 = 202
$ws.Range("E203").Value2 = 45277.04166666666
$ws.Range("F203").Value2 = 'Rosario Central'
$ws.Range("G203").Value2 = 1
$ws.Range("H203").Value2 = 'Platense'
$ws.Range("I203").Value2 = 0
$ws.Range("J203").Value2 = 3.23
$ws.Range("K203").Value2 = '11/12/2023 19:11'
$ws.Range("L203").Value2 = 2.74
$ws.Range("M203").Value2 = '17/12/2023 00:59'
$ws.Range("N203").Value2 = 2.89
$ws.Range("O203").Value2 = '11/12/2023 19:11'
$ws.Range("P203").Value2 = 2.8
$ws.Range("Q203").Value2 = '17/12/2023 00:48'
$ws.Range("R203").Value2 = 2.45
$ws.Range("S203").Value2 = '11/12/2023 19:11'
$ws.Range("T203").Value2 = 3.21
$ws.Range("U203").Value2 = '17/12/2023 00:59'
$ws.Range("V203").Value2 = 'https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/rosario-central-platense/65Pyeyh3/'
